$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 4-9 with higher precision values
$data = @{
    4  = @(12.0990725459, 44.4889279027, 11.1841151296, 8.357796220099999, 134.4981458801, 26.204610738, 97.817770129, 9.459943475299999, 10.6058217383, 1.318990554)
    5  = @(11.973591, 48.287198, 11.442368, 9.128698, 132.817807, 26.095552, 100.086072, 9.660145, 10.611959, 1.331172)
    6  = @(12.1310620997, 50.1257575583, 11.4591605393, 8.8901179019, 130.1200912182, 26.6678997134, 99.1885421109, 10.0149840289, 10.4177313761, 1.3312070994)
    7  = @(12.2222706985, 52.0905758601, 11.3856864905, 9.0029642792, 127.2272746739, 29.5234904034, 96.0983860123, 9.7036011216, 9.643514291400001, 1.3052155812)
    8  = @(12.5072369149, 56.4233866835, 13.5687424738, 10.795968135, 130.1094332834, 26.9056655875, 98.60583190680001, 10.6606203719, 9.473508819599999, 1.2725812361)
    9  = @(12.9887044066315, 56.2596316680732, 13.8778304909106, 12.747180303307, 141.188930996034, 24.7893999499941, 103.664765283476, 12.762604742767, 10.3658582266933, 1.28588624224285)
    10 = @(14.4, 61, 14.2, 12.3, 144.6, 32.9, 109.8, 13.2, 10.8, 1.3)
}

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i  # column B = 2
        $ws.Cells.Item($rowNum, $col).Value = $vals[$i]
    }
}

# Add the new label for row 10, mirroring the style of the existing label cells (A4:A9)
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A10").Value = "2021年"
